$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "306.66"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "-2.87%"
$ws.Cells.Item(2, 6).NumberFormat = "@"
$ws.Cells.Item(2, 6).Value = "17-2-2023"
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "1"

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "48.53"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "7.77%"
$ws.Cells.Item(3, 6).NumberFormat = "@"
$ws.Cells.Item(3, 6).Value = "17-2-2023"
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "1"

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "5.105"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "-1.53%"
$ws.Cells.Item(4, 6).NumberFormat = "@"
$ws.Cells.Item(4, 6).Value = "17-2-2023"
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = "1"

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.07699"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "-4.56%"
$ws.Cells.Item(5, 6).NumberFormat = "@"
$ws.Cells.Item(5, 6).Value = "17-2-2023"
$ws.Cells.Item(5, 7).NumberFormat = "@"
$ws.Cells.Item(5, 7).Value = "1"

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "4.416"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "-3.23%"
$ws.Cells.Item(6, 6).NumberFormat = "@"
$ws.Cells.Item(6, 6).Value = "17-2-2023"
$ws.Cells.Item(6, 7).NumberFormat = "@"
$ws.Cells.Item(6, 7).Value = "1"

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.307"
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "19.83%"
$ws.Cells.Item(7, 6).NumberFormat = "@"
$ws.Cells.Item(7, 6).Value = "17-2-2023"
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Cells.Item(7, 7).Value = "1"

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "1.563"
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "-6.63%"
$ws.Cells.Item(8, 6).NumberFormat = "@"
$ws.Cells.Item(8, 6).Value = "17-2-2023"
$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = "1"

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.1228"
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "-7.12%"
$ws.Cells.Item(9, 6).NumberFormat = "@"
$ws.Cells.Item(9, 6).Value = "17-2-2023"
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = "1"

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.1921"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "-0.81%"
$ws.Cells.Item(10, 6).NumberFormat = "@"
$ws.Cells.Item(10, 6).Value = "17-2-2023"
$ws.Cells.Item(10, 7).NumberFormat = "@"
$ws.Cells.Item(10, 7).Value = "1"

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.09212"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "-2.39%"
$ws.Cells.Item(11, 6).NumberFormat = "@"
$ws.Cells.Item(11, 6).Value = "17-2-2023"
$ws.Cells.Item(11, 7).NumberFormat = "@"
$ws.Cells.Item(11, 7).Value = "1"

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.04492"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "3.22%"
$ws.Cells.Item(12, 6).NumberFormat = "@"
$ws.Cells.Item(12, 6).Value = "17-2-2023"
$ws.Cells.Item(12, 7).NumberFormat = "@"
$ws.Cells.Item(12, 7).Value = "1"

$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "0.92%"
$ws.Cells.Item(13, 6).NumberFormat = "@"
$ws.Cells.Item(13, 6).Value = "17-2-2023"
$ws.Cells.Item(13, 7).NumberFormat = "@"
$ws.Cells.Item(13, 7).Value = "1"

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.001271"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "-3.27%"
$ws.Cells.Item(14, 6).NumberFormat = "@"
$ws.Cells.Item(14, 6).Value = "17-2-2023"
$ws.Cells.Item(14, 7).NumberFormat = "@"
$ws.Cells.Item(14, 7).Value = "1"

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.04168"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "-2.84%"
$ws.Cells.Item(15, 6).NumberFormat = "@"
$ws.Cells.Item(15, 6).Value = "17-2-2023"
$ws.Cells.Item(15, 7).NumberFormat = "@"
$ws.Cells.Item(15, 7).Value = "1"

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.005858"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "1.24%"
$ws.Cells.Item(16, 6).NumberFormat = "@"
$ws.Cells.Item(16, 6).Value = "17-2-2023"
$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = "1"

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "3.354"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "-2.01%"
$ws.Cells.Item(17, 6).NumberFormat = "@"
$ws.Cells.Item(17, 6).Value = "17-2-2023"
$ws.Cells.Item(17, 7).NumberFormat = "@"
$ws.Cells.Item(17, 7).Value = "1"

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "2.248"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "-7.32%"
$ws.Cells.Item(18, 6).NumberFormat = "@"
$ws.Cells.Item(18, 6).Value = "17-2-2023"
$ws.Cells.Item(18, 7).NumberFormat = "@"
$ws.Cells.Item(18, 7).Value = "1"

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.3493"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "3.32%"
$ws.Cells.Item(19, 6).NumberFormat = "@"
$ws.Cells.Item(19, 6).Value = "17-2-2023"
$ws.Cells.Item(19, 7).NumberFormat = "@"
$ws.Cells.Item(19, 7).Value = "1"

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "8.123"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "-1.53%"
$ws.Cells.Item(20, 6).NumberFormat = "@"
$ws.Cells.Item(20, 6).Value = "17-2-2023"
$ws.Cells.Item(20, 7).NumberFormat = "@"
$ws.Cells.Item(20, 7).Value = "1"

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.1383"
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "1.15%"
$ws.Cells.Item(21, 6).NumberFormat = "@"
$ws.Cells.Item(21, 6).Value = "17-2-2023"
$ws.Cells.Item(21, 7).NumberFormat = "@"
$ws.Cells.Item(21, 7).Value = "1"

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.3013"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "-4.23%"
$ws.Cells.Item(22, 6).NumberFormat = "@"
$ws.Cells.Item(22, 6).Value = "17-2-2023"
$ws.Cells.Item(22, 7).NumberFormat = "@"
$ws.Cells.Item(22, 7).Value = "1"

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.001269"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "-1.98%"
$ws.Cells.Item(23, 6).NumberFormat = "@"
$ws.Cells.Item(23, 6).Value = "17-2-2023"
$ws.Cells.Item(23, 7).NumberFormat = "@"
$ws.Cells.Item(23, 7).Value = "1"

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.004117"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "-2.70%"
$ws.Cells.Item(24, 6).NumberFormat = "@"
$ws.Cells.Item(24, 6).Value = "17-2-2023"
$ws.Cells.Item(24, 7).NumberFormat = "@"
$ws.Cells.Item(24, 7).Value = "1"

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.0001352"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "0.50%"
$ws.Cells.Item(25, 6).NumberFormat = "@"
$ws.Cells.Item(25, 6).Value = "17-2-2023"
$ws.Cells.Item(25, 7).NumberFormat = "@"
$ws.Cells.Item(25, 7).Value = "1"

$ws.Cells.Item(26, 6).NumberFormat = "@"
$ws.Cells.Item(26, 6).Value = "17-2-2023"
$ws.Cells.Item(26, 7).NumberFormat = "@"
$ws.Cells.Item(26, 7).Value = "1"

$ws.Cells.Item(27, 6).NumberFormat = "@"
$ws.Cells.Item(27, 6).Value = "17-2-2023"
$ws.Cells.Item(27, 7).NumberFormat = "@"
$ws.Cells.Item(27, 7).Value = "1"

$ws.Cells.Item(28, 6).NumberFormat = "@"
$ws.Cells.Item(28, 6).Value = "17-2-2023"
$ws.Cells.Item(28, 7).NumberFormat = "@"
$ws.Cells.Item(28, 7).Value = "1"

$ws.Cells.Item(29, 6).NumberFormat = "@"
$ws.Cells.Item(29, 6).Value = "17-2-2023"
$ws.Cells.Item(29, 7).NumberFormat = "@"
$ws.Cells.Item(29, 7).Value = "1"

$ws.Cells.Item(30, 6).NumberFormat = "@"
$ws.Cells.Item(30, 6).Value = "17-2-2023"
$ws.Cells.Item(30, 7).NumberFormat = "@"
$ws.Cells.Item(30, 7).Value = "1"

$ws.Cells.Item(31, 6).NumberFormat = "@"
$ws.Cells.Item(31, 6).Value = "17-2-2023"
$ws.Cells.Item(31, 7).NumberFormat = "@"
$ws.Cells.Item(31, 7).Value = "1"

$ws.Cells.Item(32, 6).NumberFormat = "@"
$ws.Cells.Item(32, 6).Value = "17-2-2023"
$ws.Cells.Item(32, 7).NumberFormat = "@"
$ws.Cells.Item(32, 7).Value = "1"

$ws.Cells.Item(33, 6).NumberFormat = "@"
$ws.Cells.Item(33, 6).Value = "17-2-2023"
$ws.Cells.Item(33, 7).NumberFormat = "@"
$ws.Cells.Item(33, 7).Value = "1"

$ws.Cells.Item(34, 6).NumberFormat = "@"
$ws.Cells.Item(34, 6).Value = "17-2-2023"
$ws.Cells.Item(34, 7).NumberFormat = "@"
$ws.Cells.Item(34, 7).Value = "1"

$ws.Cells.Item(35, 6).NumberFormat = "@"
$ws.Cells.Item(35, 6).Value = "17-2-2023"
$ws.Cells.Item(35, 7).NumberFormat = "@"
$ws.Cells.Item(35, 7).Value = "1"

$ws.Cells.Item(36, 6).NumberFormat = "@"
$ws.Cells.Item(36, 6).Value = "17-2-2023"
$ws.Cells.Item(36, 7).NumberFormat = "@"
$ws.Cells.Item(36, 7).Value = "1"

$ws.Cells.Item(37, 6).NumberFormat = "@"
$ws.Cells.Item(37, 6).Value = "17-2-2023"
$ws.Cells.Item(37, 7).NumberFormat = "@"
$ws.Cells.Item(37, 7).Value = "1"

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.02557"
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = "-4.37%"
$ws.Cells.Item(38, 6).NumberFormat = "@"
$ws.Cells.Item(38, 6).Value = "17-2-2023"
$ws.Cells.Item(38, 7).NumberFormat = "@"
$ws.Cells.Item(38, 7).Value = "1"

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.05714"
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "5.08%"
$ws.Cells.Item(39, 6).NumberFormat = "@"
$ws.Cells.Item(39, 6).Value = "17-2-2023"
$ws.Cells.Item(39, 7).NumberFormat = "@"
$ws.Cells.Item(39, 7).Value = "1"

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.01078"
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "83.85%"
$ws.Cells.Item(40, 6).NumberFormat = "@"
$ws.Cells.Item(40, 6).Value = "17-2-2023"
$ws.Cells.Item(40, 7).NumberFormat = "@"
$ws.Cells.Item(40, 7).Value = "1"

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.007924"
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "2.60%"
$ws.Cells.Item(41, 6).NumberFormat = "@"
$ws.Cells.Item(41, 6).Value = "17-2-2023"
$ws.Cells.Item(41, 7).NumberFormat = "@"
$ws.Cells.Item(41, 7).Value = "1"

$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "-1.52%"
$ws.Cells.Item(42, 6).NumberFormat = "@"
$ws.Cells.Item(42, 6).Value = "17-2-2023"
$ws.Cells.Item(42, 7).NumberFormat = "@"
$ws.Cells.Item(42, 7).Value = "1"

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.008458"
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "15.03%"
$ws.Cells.Item(43, 6).NumberFormat = "@"
$ws.Cells.Item(43, 6).Value = "17-2-2023"
$ws.Cells.Item(43, 7).NumberFormat = "@"
$ws.Cells.Item(43, 7).Value = "1"

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.007719"
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "-9.95%"
$ws.Cells.Item(44, 6).NumberFormat = "@"
$ws.Cells.Item(44, 6).Value = "17-2-2023"
$ws.Cells.Item(44, 7).NumberFormat = "@"
$ws.Cells.Item(44, 7).Value = "1"

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.3333"
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "6.14%"
$ws.Cells.Item(45, 6).NumberFormat = "@"
$ws.Cells.Item(45, 6).Value = "17-2-2023"
$ws.Cells.Item(45, 7).NumberFormat = "@"
$ws.Cells.Item(45, 7).Value = "1"

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.00006800"
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "-1.09%"
$ws.Cells.Item(46, 6).NumberFormat = "@"
$ws.Cells.Item(46, 6).Value = "17-2-2023"
$ws.Cells.Item(46, 7).NumberFormat = "@"
$ws.Cells.Item(46, 7).Value = "1"

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.00000000751"
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "0.53%"
$ws.Cells.Item(47, 6).NumberFormat = "@"
$ws.Cells.Item(47, 6).Value = "17-2-2023"
$ws.Cells.Item(47, 7).NumberFormat = "@"
$ws.Cells.Item(47, 7).Value = "1"

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.05666"
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = "-7.53%"
$ws.Cells.Item(48, 6).NumberFormat = "@"
$ws.Cells.Item(48, 6).Value = "17-2-2023"
$ws.Cells.Item(48, 7).NumberFormat = "@"
$ws.Cells.Item(48, 7).Value = "1"

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.004005"
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = "0.49%"
$ws.Cells.Item(49, 6).NumberFormat = "@"
$ws.Cells.Item(49, 6).Value = "17-2-2023"
$ws.Cells.Item(49, 7).NumberFormat = "@"
$ws.Cells.Item(49, 7).Value = "1"

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.00002104"
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = "0.53%"
$ws.Cells.Item(50, 6).NumberFormat = "@"
$ws.Cells.Item(50, 6).Value = "17-2-2023"
$ws.Cells.Item(50, 7).NumberFormat = "@"
$ws.Cells.Item(50, 7).Value = "1"

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0002003"
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = "0.53%"
$ws.Cells.Item(51, 6).NumberFormat = "@"
$ws.Cells.Item(51, 6).Value = "17-2-2023"
$ws.Cells.Item(51, 7).NumberFormat = "@"
$ws.Cells.Item(51, 7).Value = "1"
